$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "64.376.47"
Set-TextValue $ws.Range("E2") "  -0.16%  "
Set-TextValue $ws.Range("D3") "3.402.83"
Set-TextValue $ws.Range("E3") "  -1.61%  "
Set-TextValue $ws.Range("E4") "  -0.01%  "
Set-TextValue $ws.Range("D5") "567.62"
Set-TextValue $ws.Range("E5") "  -1.05%  "
Set-TextValue $ws.Range("D6") "156.01"
Set-TextValue $ws.Range("E6") "  -2.32%  "
Set-TextValue $ws.Range("E7") "  +0.12%  "
Set-TextValue $ws.Range("D8") "0.606"
Set-TextValue $ws.Range("E8") "  +3.66%  "
Set-TextValue $ws.Range("D9") "3.401.62"
Set-TextValue $ws.Range("E9") "  -1.88%  "
Set-TextValue $ws.Range("D10") "7.20"
Set-TextValue $ws.Range("E10") "  -2.20%  "
Set-TextValue $ws.Range("E11") "  -2.75%  "
Set-TextValue $ws.Range("D12") "0.442"
Set-TextValue $ws.Range("E12") "  -0.95%  "
Set-TextValue $ws.Range("D13") "3.990.04"
Set-TextValue $ws.Range("E13") "  -1.55%  "
Set-TextValue $ws.Range("E14") "  +0.08%  "
Set-TextValue $ws.Range("D15") "0.0000187"
Set-TextValue $ws.Range("E15") "  -3.58%  "
Set-TextValue $ws.Range("D16") "27.68"
Set-TextValue $ws.Range("E16") "  -3.64%  "
Set-TextValue $ws.Range("D17") "64.448.12"
Set-TextValue $ws.Range("D18") "3.413.56"
Set-TextValue $ws.Range("E18") "  -1.60%  "
Set-TextValue $ws.Range("D19") "6.34"
Set-TextValue $ws.Range("E19") "  -1.14%  "
Set-TextValue $ws.Range("D20") "13.96"
Set-TextValue $ws.Range("E20") "  -2.93%  "
Set-TextValue $ws.Range("D21") "374.71"
Set-TextValue $ws.Range("E21") "  -3.10%  "
Set-TextValue $ws.Range("D22") "7.97"
Set-TextValue $ws.Range("E22") "  -2.51%  "
Set-TextValue $ws.Range("D23") "0.548"
Set-TextValue $ws.Range("E23") "  +0.81%  "
Set-TextValue $ws.Range("D24") "0.998"
Set-TextValue $ws.Range("E24") "  -0.31%  "
Set-TextValue $ws.Range("D25") "72.01"
Set-TextValue $ws.Range("E25") "  -1.69%  "
Set-TextValue $ws.Range("D26") "0.0000117"
Set-TextValue $ws.Range("E26") "  -2.78%  "
Set-TextValue $ws.Range("D27") "9.90"
Set-TextValue $ws.Range("E27") "  +4.04%  "
Set-TextValue $ws.Range("D28") "0.175"
Set-TextValue $ws.Range("E28") "  -2.36%  "
Set-TextValue $ws.Range("E29") "  -0.20%  "
Set-TextValue $ws.Range("D30") "1.46"
Set-TextValue $ws.Range("E30") "  +1.34%  "
Set-TextValue $ws.Range("D31") "6.06"
Set-TextValue $ws.Range("E31") "  -1.88%  "
Set-TextValue $ws.Range("D32") "2.02"
Set-TextValue $ws.Range("E32") "  -0.42%  "
Set-TextValue $ws.Range("D33") "23.16"
Set-TextValue $ws.Range("E33") "  -2.04%  "
Set-TextValue $ws.Range("D34") "7.14"
Set-TextValue $ws.Range("E34") "  +1.44%  "
Set-TextValue $ws.Range("D35") "1.58"
Set-TextValue $ws.Range("E35") "  +5.45%  "
Set-TextValue $ws.Range("D36") "159.69"
Set-TextValue $ws.Range("E36") "  -0.65%  "
Set-TextValue $ws.Range("D37") "1.88"
Set-TextValue $ws.Range("E37") "  -0.10%  "
Set-TextValue $ws.Range("D38") "0.0759"
Set-TextValue $ws.Range("E38") "  -2.09%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D39") "6.80"
Set-TextValue $ws.Range("E39") "  +2.04%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D40") "26.78"
Set-TextValue $ws.Range("E40") "  -1.75%  "
Set-TextValue $ws.Range("D41") "4.66"
Set-TextValue $ws.Range("E41") "  +4.48%  "
Set-TextValue $ws.Range("D42") "2.830.89"
Set-TextValue $ws.Range("E42") "  -2.75%  "
Set-TextValue $ws.Range("D43") "42.70"
Set-TextValue $ws.Range("E43") "  +1.26%  "
Set-TextValue $ws.Range("D44") "0.0313"
Set-TextValue $ws.Range("E44") "  -1.85%  "
Set-TextValue $ws.Range("D45") "0.765"
Set-TextValue $ws.Range("E45") "  -0.92%  "
Set-TextValue $ws.Range("D46") "25.59"
Set-TextValue $ws.Range("E46") "  +7.56%  "
Set-TextValue $ws.Range("D47") "1.07"
Set-TextValue $ws.Range("E47") "  -2.46%  "
Set-TextValue $ws.Range("D48") "309.08"
Set-TextValue $ws.Range("E48") "  +4.18%  "
Set-TextValue $ws.Range("D49") "0.108"
Set-TextValue $ws.Range("E49") "  +0.23%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D50") "6.54"
Set-TextValue $ws.Range("E50") "  +0.22%  "
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue $ws.Range("D51") "0.857"
Set-TextValue $ws.Range("E51") "  -0.75%  "
